# projectTimetable.xlsx update:
#  - Task 2 / Task 3 descriptions gain lead/helper attributions
#  - B4 picks up the same "in progress" green fill as C4:E4 (Task 2 now
#    starts a week earlier now that Andrew Samuel leads it)
#  - Selection moves to A5 (Task 3 row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task text updates (shared strings) ---
$ws.Range("A4").Value = "Task 2: Write code using SQL for the database (Andrew Samuel (lead), Ary Hernandez, Jacquelyn Johnson)"
$ws.Range("A5").Value = "Task 3: Implement database within an application (Jacquelyn Johnson (lead), Andrew Samuel, Ary Hernandez) "

# --- B4 fill: match the highlighted style already used by C4:E4 ---
$ws.Range("C4").Copy()
$ws.Range("B4").PasteSpecial(-4122)   # xlPasteFormats

# --- Selection moves to A5 ---
$ws.Range("A5").Select() | Out-Null
